$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Text-formatted columns (B,C,D,E) are written with NumberFormat "@"
# first so numeric-looking strings (e.g. "1.00", "89.828.26") are not
# coerced into numbers by Excel's type inference on .Value assignment.

$cell = $ws.Range("D2")
$cell.NumberFormat = '@'
$cell.Value = '90.336.95'

$cell = $ws.Range("E2")
$cell.NumberFormat = '@'
$cell.Value = '  +0.13%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = '@'
$cell.Value = '3.102.24'

$cell = $ws.Range("E3")
$cell.NumberFormat = '@'
$cell.Value = '  +0.65%  '

$cell = $ws.Range("E4")
$cell.NumberFormat = '@'
$cell.Value = '  -0.15%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = '@'
$cell.Value = '240.79'

$cell = $ws.Range("E5")
$cell.NumberFormat = '@'
$cell.Value = '  +3.79%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = '@'
$cell.Value = '619.90'

$cell = $ws.Range("E6")
$cell.NumberFormat = '@'
$cell.Value = '  -0.60%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = '@'
$cell.Value = '1.14'

$cell = $ws.Range("E7")
$cell.NumberFormat = '@'
$cell.Value = '  +3.77%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = '@'
$cell.Value = '0.364'

$cell = $ws.Range("E8")
$cell.NumberFormat = '@'
$cell.Value = '  +1.10%  '

$cell = $ws.Range("E9")
$cell.NumberFormat = '@'
$cell.Value = '  -0.09%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = '@'
$cell.Value = '3.099.36'

$cell = $ws.Range("E10")
$cell.NumberFormat = '@'
$cell.Value = '  +24.07%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = '@'
$cell.Value = '0.734'

$cell = $ws.Range("E11")
$cell.NumberFormat = '@'
$cell.Value = '  +0.89%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = '@'
$cell.Value = '0.203'

$cell = $ws.Range("E12")
$cell.NumberFormat = '@'
$cell.Value = '  +3.40%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = '@'
$cell.Value = '0.0000247'

$cell = $ws.Range("E13")
$cell.NumberFormat = '@'
$cell.Value = '  +0.14%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = '@'
$cell.Value = '34.93'

$cell = $ws.Range("E14")
$cell.NumberFormat = '@'
$cell.Value = '  -3.71%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = '@'
$cell.Value = '5.49'

$cell = $ws.Range("E15")
$cell.NumberFormat = '@'
$cell.Value = '  +0.52%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = '@'
$cell.Value = '90.177.00'

$cell = $ws.Range("E16")
$cell.NumberFormat = '@'
$cell.Value = '  +0.06%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = '@'
$cell.Value = '3.670.50'

$cell = $ws.Range("E17")
$cell.NumberFormat = '@'
$cell.Value = '  +0.31%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = '@'
$cell.Value = '3.104.14'

$cell = $ws.Range("E18")
$cell.NumberFormat = '@'
$cell.Value = '  +0.58%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = '@'
$cell.Value = '3.68'

$cell = $ws.Range("E19")
$cell.NumberFormat = '@'
$cell.Value = '  -1.98%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = '@'
$cell.Value = '14.69'

$cell = $ws.Range("E20")
$cell.NumberFormat = '@'
$cell.Value = '  +4.98%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = '@'
$cell.Value = '0.0000210'

$cell = $ws.Range("E21")
$cell.NumberFormat = '@'
$cell.Value = '  +0.69%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = '@'
$cell.Value = '5.81'

$cell = $ws.Range("E22")
$cell.NumberFormat = '@'
$cell.Value = '  +4.73%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = '@'
$cell.Value = '438.38'

$cell = $ws.Range("E23")
$cell.NumberFormat = '@'
$cell.Value = '  +0.11%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = '@'
$cell.Value = '9.07'

$cell = $ws.Range("E24")
$cell.NumberFormat = '@'
$cell.Value = '  +2.29%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = '@'
$cell.Value = '5.50'

$cell = $ws.Range("E25")
$cell.NumberFormat = '@'
$cell.Value = '  -3.21%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = '@'
$cell.Value = '89.83'

$cell = $ws.Range("E26")
$cell.NumberFormat = '@'
$cell.Value = '  +1.02%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = '@'
$cell.Value = '11.87'

$cell = $ws.Range("E27")
$cell.NumberFormat = '@'
$cell.Value = '  -2.37%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = '@'
$cell.Value = '3.227.99'

$cell = $ws.Range("D29")
$cell.NumberFormat = '@'
$cell.Value = '1.00'

$cell = $ws.Range("E29")
$cell.NumberFormat = '@'
$cell.Value = '  +0.02%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = '@'
$cell.Value = '0.245'

$cell = $ws.Range("E30")
$cell.NumberFormat = '@'
$cell.Value = '  +21.26%  '

$cell = $ws.Range("E31")
$cell.NumberFormat = '@'
$cell.Value = '  +10.78%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = '@'
$cell.Value = '0.120'

$cell = $ws.Range("E32")
$cell.NumberFormat = '@'
$cell.Value = '  +32.41%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = '@'
$cell.Value = '9.27'

$cell = $ws.Range("E33")
$cell.NumberFormat = '@'
$cell.Value = '  -1.47%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = '@'
$cell.Value = '0.166'

$cell = $ws.Range("E34")
$cell.NumberFormat = '@'
$cell.Value = '  +8.25%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = '@'
$cell.Value = '0.967'

$cell = $ws.Range("E35")
$cell.NumberFormat = '@'
$cell.Value = '  -0.99%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = '@'
$cell.Value = '7.90'

$cell = $ws.Range("E36")
$cell.NumberFormat = '@'
$cell.Value = '  +13.12%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = '@'
$cell.Value = '4.28'

$cell = $ws.Range("E37")
$cell.NumberFormat = '@'
$cell.Value = '  +20.91%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = '@'
$cell.Value = '26.24'

$cell = $ws.Range("E38")
$cell.NumberFormat = '@'
$cell.Value = '  +0.29%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = '@'
$cell.Value = '1.92'

$cell = $ws.Range("E39")
$cell.NumberFormat = '@'
$cell.Value = '  +0.48%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = '@'
$cell.Value = '491.13'

$cell = $ws.Range("E40")
$cell.NumberFormat = '@'
$cell.Value = '  -2.67%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = '@'
$cell.Value = '3.56'

$cell = $ws.Range("E41")
$cell.NumberFormat = '@'
$cell.Value = '  -5.30%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = '@'
$cell.Value = '1.30'

$cell = $ws.Range("E42")
$cell.NumberFormat = '@'
$cell.Value = '  +2.34%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = '@'
$cell.Value = '0.420'

$cell = $ws.Range("E43")
$cell.NumberFormat = '@'
$cell.Value = '  +2.68%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = '@'
$cell.Value = '22.15'

$cell = $ws.Range("E44")
$cell.NumberFormat = '@'
$cell.Value = '  -0.12%  '

$cell = $ws.Range("E45")
$cell.NumberFormat = '@'
$cell.Value = '  +0.01%  '

$cell = $ws.Range("B46")
$cell.NumberFormat = '@'
$cell.Value = 'Stacks'

$cell = $ws.Range("C46")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'

$cell = $ws.Range("D46")
$cell.NumberFormat = '@'
$cell.Value = '1.93'

$cell = $ws.Range("E46")
$cell.NumberFormat = '@'
$cell.Value = '  +1.64%  '

$cell = $ws.Range("B47")
$cell.NumberFormat = '@'
$cell.Value = 'Monero'

$cell = $ws.Range("C47")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$cell = $ws.Range("D47")
$cell.NumberFormat = '@'
$cell.Value = '155.19'

$cell = $ws.Range("E47")
$cell.NumberFormat = '@'
$cell.Value = '  +2.96%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = '@'
$cell.Value = '0.691'

$cell = $ws.Range("E48")
$cell.NumberFormat = '@'
$cell.Value = '  +0.84%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = '@'
$cell.Value = '1.35'

$cell = $ws.Range("E49")
$cell.NumberFormat = '@'
$cell.Value = '  +1.32%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = '@'
$cell.Value = '44.22'

$cell = $ws.Range("E50")
$cell.NumberFormat = '@'
$cell.Value = '  -1.64%  '

$cell = $ws.Range("B51")
$cell.NumberFormat = '@'
$cell.Value = 'Filecoin'

$cell = $ws.Range("C51")
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'

$cell = $ws.Range("D51")
$cell.NumberFormat = '@'
$cell.Value = '4.41'

$cell = $ws.Range("E51")
$cell.NumberFormat = '@'
$cell.Value = '  -0.39%  '
